$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 83: EMI Filter 100pF, no ferrite ---
# D83 remark text gets a clarification appended
$ws.Range("D83").Value = "Murata DSN6NC51H101, replacement is with ferrite"

# F83 is a brand-new cell (Mouser part link) that did not exist before
$ws.Range("F83").Value = "81-DSS1NB32A101Q91A"
$ws.Hyperlinks.Add($ws.Range("F83"), "https://www.mouser.com/ProductDetail/81-DSS1NB32A101Q91A")
$ws.Range("F83").Style = "Hyperlink"

# --- Row 90: Ferrite Bead ---
# F90 is a brand-new cell (Mouser part link) that did not exist before
$ws.Range("F90").Value = "542-FB20020-4B-RC"
$ws.Hyperlinks.Add($ws.Range("F90"), "https://www.mouser.com/ProductDetail/542-FB20020-4B-RC")
$ws.Range("F90").Style = "Hyperlink"

# --- Row 91: Ferrite Bead Long ---
# F91 Mouser part number changes, update text + hyperlink target
$ws.Range("F91").Value = "434-BEAD-10-600P-02"
$ws.Hyperlinks.Add($ws.Range("F91"), "https://www.mouser.com/ProductDetail/434-BEAD-10-600P-02")
$ws.Range("F91").Style = "Hyperlink"

# --- Row 103: 5 Pin Square DIN (CN8) ---
# D103 remark gets a note about the Ramixx500 alternative below
$ws.Range("D103").Value = "Power Dynamics DS-215. Rämixx500 see below."

# --- New row 141: CN8 alternative connector (6 Pin Round DIN) ---
$ws.Rows.Item(141).Insert()
$ws.Range("B141").Value = 1
$ws.Range("C141").Value = "6 Pin Round DIN"
$ws.Range("D141").Value = "CN8 alternative. Requires modification to PSU cable!"
$ws.Range("F141").Value = "490-SDF-60J"
$ws.Hyperlinks.Add($ws.Range("F141"), "https://www.mouser.com/ProductDetail/490-SDF-60J")
$ws.Range("F141").Style = "Hyperlink"
$ws.Range("G141").Value = "CN8"
$ws.Range("G141").WrapText = $true
